$wb = $excel.ActiveWorkbook

# --- Sheet3 (CheckInVisitors): drop the "active/scrolled" view state ---
# (Adding the new sheets below and making the last one active will naturally
# clear tabSelected/topLeftCell from this sheet's view.)

# --- New sheet: CheckedIn ---
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$wsCheckedIn = $wb.Worksheets.Add($null, $lastSheet)
$wsCheckedIn.Name = "CheckedIn"
$wsCheckedIn.Columns.Item(1).ColumnWidth = 12.944010416666666
$wsCheckedIn.Range("A1").Value = "visitor_NIC"
$wsCheckedIn.Range("A2").Value = "6348445764v"
$wsCheckedIn.Range("A3").Value = "7821459632v"
$wsCheckedIn.Range("A4").Value = "4578963245v"
$wsCheckedIn.Range("A1").Interior.Color = 65535

# --- New sheet: Overdue ---
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$wsOverdue = $wb.Worksheets.Add($null, $lastSheet)
$wsOverdue.Name = "Overdue"
$wsOverdue.Columns.Item(1).ColumnWidth = 11.944010416666666
$wsOverdue.Range("A1").Value = "visitor_NIC"
$wsOverdue.Range("A2").Value = "5612345783v"
$wsOverdue.Range("A3").Value = "8695748612v"
$wsOverdue.Range("A1").Interior.Color = 65535

# --- New sheet: manageBuilding ---
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$wsBuilding = $wb.Worksheets.Add($null, $lastSheet)
$wsBuilding.Name = "manageBuilding"
$wsBuilding.Columns.Item(1).ColumnWidth = 18.276041666666668
$wsBuilding.Columns.Item(2).ColumnWidth = 14.166666666666666
$wsBuilding.Range("A1").Value = "BuildingName"
$wsBuilding.Range("B1").Value = "floorName"
$wsBuilding.Range("A2").Value = "buildingNo1"
$wsBuilding.Range("A3").Value = "buildingNo2"
$wsBuilding.Range("B2").Value = "floorNo1"
$wsBuilding.Range("B3").Value = "floorNo2"
$wsBuilding.Range("A1:B1").Interior.Color = 65535

Write-Output "done"
